# TC32_Canine_StudyNCATSCOP01-StudyPart_SampleType_FileType.xlsx
# "NCATS Study multifilter testcases61to70"
#
# The author re-entered the "Sample ID" Neo4j query text into cell B3 on the
# "startup" sheet (content unchanged) and left the selection on B3 when the
# file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the (unchanged) Sample ID query text to B3 so it is freshly
# (re)written to the sheet, exactly as it was before.
$sampleIdQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`n OPTIONAL MATCH (c)-->(ci:canine_individual)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p,ci`n  where s.clinical_study_designation IN ['NCATS-COP01'] and ci IS NULL and samp.summarized_sample_type in ['Whole Blood'] and sf.file_type in ['Supplemental Data File']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`noptional MATCH (samp:sample)-->(c) `nWITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed,`n        coalesce(diag.disease_term,'') AS Diagnosis, `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation```norder by samp.sample_id asc`nlimit 100"

# Remember the current row height - re-writing the cell value can trigger an
# autofit recalculation, and the row height is not meant to change here.
$row3Height = $ws.Rows.Item(3).RowHeight

$ws.Range("B3").Value = $sampleIdQuery

# Restore row 3's height in case setting the value above auto-adjusted it.
$ws.Rows.Item(3).RowHeight = $row3Height

# Leave the selection on B3, matching the saved cursor position.
$ws.Range("B3").Select()
